$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Multivalued" column (K) ---

# Header cell K4: bold header style matching the other header cells (row 4),
# but with black color + Calibri "swiss" family, like the workbook title font.
$ws.Range("K4").Value = "Multivalued"
$ws.Range("A5").Copy()
$ws.Range("K4").PasteSpecial(-4122)
$ws.Range("K4").Font.Bold = $true
$ws.Range("K4").Font.Color = 0
$ws.Range("K4").Font.Family = 2

# Data cells K5:K8: literal "FALSE" text (kept as text, not boolean),
# formatted with a TRUE/FALSE custom display format, left aligned.
$ws.Range("K5:K8").Value = "'FALSE"
$ws.Range("K5:K8").NumberFormat = """TRUE"";""TRUE"";""FALSE"""
$ws.Range("K5:K8").HorizontalAlignment = -4131

# --- Selection / view state ---
$ws.Range("L8").Select()
